$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (B1:K1) -- A1 is already blank in the source and is left untouched
$ws.Range("B1").Value = '$ bold(''All'')'
$ws.Range("C1").Value = 'Millionaires'
$ws.Range("D1").Value = 'Japan Non-voters'
$ws.Range("E1").Value = 'Japan Left'
$ws.Range("F1").Value = 'Japan Center/Right'
$ws.Range("G1").Value = 'Saudi Arabia'
$ws.Range("H1").Value = 'Saudi citizens'
$ws.Range("I1").Value = 'U.S. Non-voters'
$ws.Range("J1").Value = 'U.S. Harris'
$ws.Range("K1").Value = 'U.S. Trump'

# Rows 2-11: policy label (A) + country/group values (B:K)
$data = New-Object 'object[,]' 10,11
$data[0,0] = 'Minimum tax of 2% on billionaires''
wealth, in voluntary countries'
$data[0,1] = 0.746758615571303
$data[0,2] = 0.36193361088839
$data[0,3] = 0.507672921271455
$data[0,4] = 0.707476508293962
$data[0,5] = 0.495218884354337
$data[0,6] = 0.844702008923008
$data[0,7] = 0.847665991116618
$data[0,8] = 0.637831843511571
$data[0,9] = 1.27187155660457
$data[0,10] = 0.0650151157240103
$data[1,0] = 'Bridgetown initiative: MDBs expanding sustainable
investments in LICs, and at lower interest rates'
$data[1,1] = 0.53880897005617
$data[1,2] = 0.490786812397021
$data[1,3] = 0.270227396821779
$data[1,4] = 0.483958742029262
$data[1,5] = 0.402985842257098
$data[1,6] = 0.834571283823028
$data[1,7] = 0.854245566734256
$data[1,8] = 0.420572659025095
$data[1,9] = 0.967666839632799
$data[1,10] = -0.0924690322664047
$data[2,0] = 'L&D: Developed countries financing a fund to help
vulnerable countries cope with climate Loss and damage'
$data[2,1] = 0.457646793262474
$data[2,2] = 0.220852126457596
$data[2,3] = 0.190001288040726
$data[2,4] = 0.369650389914943
$data[2,5] = 0.315628598381185
$data[2,6] = 0.958227224438065
$data[2,7] = 0.979381082032853
$data[2,8] = 0.3751153596592
$data[2,9] = 0.860594629094862
$data[2,10] = -0.224011634614708
$data[3,0] = 'At least 0.7% of developed countries'' GDP in foreign aid'
$data[3,1] = 0.339613614975512
$data[3,2] = 0.202877418099734
$data[3,3] = 0.0130537361393236
$data[3,4] = 0.170380103753737
$data[3,5] = 0.142327045269292
$data[3,6] = 0.808599687276629
$data[3,7] = 0.840369921573375
$data[3,8] = 0.217943990590892
$data[3,9] = 0.711174261698852
$data[3,10] = -0.199872640439334
$data[4,0] = 'Debt relief for vulnerable countries, suspending
payments until they are more able to repay'
$data[4,1] = 0.338691995923963
$data[4,2] = -0.00161372574507126
$data[4,3] = 0.156007282015151
$data[4,4] = 0.340517416281124
$data[4,5] = 0.178406875332954
$data[4,6] = 0.932399297563257
$data[4,7] = 0.983026288814376
$data[4,8] = 0.359431660781135
$data[4,9] = 0.635303602175847
$data[4,10] = -0.173015941873683
$data[5,0] = 'Raise global minimum tax on profit from 15% to 35%,
allocating revenues to countries based on sales'
$data[5,1] = 0.329110709214005
$data[5,2] = 0.149439461021209
$data[5,3] = 0.215238393472086
$data[5,4] = 0.475226190923084
$data[5,5] = 0.314508297057087
$data[5,6] = 0.497163901847429
$data[5,7] = 0.544096133581187
$data[5,8] = 0.227936339106368
$data[5,9] = 0.760894205533922
$data[5,10] = -0.258309044847189
$data[6,0] = 'NCQG: Developing countries providing $300 bn a
year in climate finance for developing countries'
$data[6,1] = 0.320103458109216
$data[6,2] = 0.0776885061155204
$data[6,3] = -0.0232948275266601
$data[6,4] = 0.167835305478371
$data[6,5] = 0.0826731315459806
$data[6,6] = 0.776081506237869
$data[6,7] = 0.799347616994926
$data[6,8] = 0.230856615873908
$data[6,9] = 0.710961677011096
$data[6,10] = -0.486893200935083
$data[7,0] = 'International levy on shipping carbon emissions,
returned to countries based on population'
$data[7,1] = 0.31788661467437
$data[7,2] = 0.12548549877193
$data[7,3] = 0.00535231156239397
$data[7,4] = 0.124308334379696
$data[7,5] = 0.0757577552841607
$data[7,6] = 0.660342805051329
$data[7,7] = 0.754267748681372
$data[7,8] = 0.226492442047215
$data[7,9] = 0.69384747595753
$data[7,10] = -0.206802852572372
$data[8,0] = 'Expand Security Council to new permanent members (e.g.
India, Brazil, African Union), restrict veto use'
$data[8,1] = 0.315567059102638
$data[8,2] = 0.17400074680061
$data[8,3] = 0.113656480813608
$data[8,4] = 0.335172796251208
$data[8,5] = 0.268968939604004
$data[8,6] = 0.707615721886186
$data[8,7] = 0.741335070219383
$data[8,8] = 0.242987426599086
$data[8,9] = 0.711947039923566
$data[8,10] = -0.221135828867629
$data[9,0] = 'International levy on aviation carbon emissions, raising
prices by 30%, returned to countries based on population'
$data[9,1] = 0.00973878914166017
$data[9,2] = -0.213557609245605
$data[9,3] = -0.121424317090429
$data[9,4] = -0.0580953959864491
$data[9,5] = -0.0790943764613775
$data[9,6] = 0.414946008655522
$data[9,7] = 0.528044529912143
$data[9,8] = -0.057232812078052
$data[9,9] = 0.328732155910792
$data[9,10] = -0.483840438600429

$ws.Range("A2:K11").Value = $data

# Remove now-unused trailing columns (L:O) that held the removed Europe data
$ws.Range("L1:O11").Clear()

Write-Output "done"